$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.367.58"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "1.847.84"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").Value = "0.9978"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "240.51"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "0.6270"
$ws.Range("E6").Value = "  -0.56%  "

$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "0.07619"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").Value = "24.78"
$ws.Range("E10").Value = "  +1.06%  "

$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "5.028"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "0.6797"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "0.00001056"
$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("D15").Value = "83.00"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "6.143"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "29.404.96"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "227.70"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").Value = "7.470"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").Value = "0.9987"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "158.89"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("D24").Value = "0.1383"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").Value = "8.431"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").Value = "17.66"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").Value = "1.413"
$ws.Range("E27").Value = "  +7.68%  "

$ws.Range("D28").Value = "1.459"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "0.05606"
$ws.Range("E29").Value = "  -2.51%  "

$ws.Range("D30").Value = "4.105"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("D31").Value = "4.069"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("B32").Value = "LidoDAOToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D32").Value = "1.835"
$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.162"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "0.6966"
$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("D35").Value = "2.587"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").Value = "1.226.68"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "2.717"
$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").Value = "0.9009"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D41").Value = "0.9993"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "101.38"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("D43").Value = "65.49"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").Value = "7.194"
$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").Value = "0.3992"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").Value = "9.006"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("E47").Value = "  -0.18%  "

$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").Value = "0.05698"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.4623"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.00000000108"
$ws.Range("E51").Value = "  -11.54%  "
